# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a 🚀
#
# Updates the StructureDefinition-employee-hire-date summary workbook:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publish timestamp
#  - replace the duplicated "Contact / No display for ContactDetail" row with
#    "Publisher: Alvearie Team" and a new "Jurisdiction: United States of America" row
#  - fix the root Extension element's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Rows 10 and 11 were a duplicated "Contact" / "No display for ContactDetail"
# pair. Delete row 11 (the duplicate) so every row below shifts up by one,
# then turn the remaining row 10 into the new "Jurisdiction" row. Doing it
# this way (delete + rename) leaves every other untouched row's cell
# type/value completely alone (e.g. the literal text "false" in the
# Abstract row stays text instead of becoming a boolean).
$meta.Rows.Item(11).Delete()

$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# --- Elements sheet --------------------------------------------------------

# Root Extension row (row 2): Short / Definition columns (K / L) describe
# this specific extension rather than the generic "Extension" text.
$elements.Cells.Item(2, 11).Value = "Employee Hire Date"
$elements.Cells.Item(2, 12).Value = "First date of employment for the employee"
